# 2023 draft big boards — add Adam Bushman's second/updated big-board sheet
# (new episode draft) and tweak the view state on the original board.

$wb = $excel.ActiveWorkbook

# Start from a duplicate of the existing "Adam-Bushman" sheet so formatting
# (column widths, date number format on column I, etc.) carries over, then
# place it immediately after the source sheet and rename it.
$src = $wb.Worksheets.Item("Adam-Bushman")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Adam-Bushman-2"

# The new board only has 44 players (rows 2:45) vs. the old 50 (rows 2:51) —
# wipe the now-unused tail rows copied over from the source sheet.
$new.Range("A46:I51").Clear()

# Write the refreshed/re-ranked big board data (pick, tier, name,
# off_position, def_position, club, height, weight, birthdate).
    $new.Cells.Item(2, 1).Value = 1; $new.Cells.Item(2, 2).Value = 1; $new.Cells.Item(2, 3).Value = 'Victor Wembanyama'; $new.Cells.Item(2, 4).Value = 'Wing'; $new.Cells.Item(2, 5).Value = 'Big'; $new.Cells.Item(2, 6).Value = 'Metropolitans92'; $new.Cells.Item(2, 7).Value = '7''4"'; $new.Cells.Item(2, 8).Value = '229 lbs'; $new.Cells.Item(2, 9).Value = 37990
    $new.Cells.Item(3, 1).Value = 2; $new.Cells.Item(3, 2).Value = 2; $new.Cells.Item(3, 3).Value = 'Brandon Miller'; $new.Cells.Item(3, 4).Value = 'Wing'; $new.Cells.Item(3, 5).Value = 'Wing'; $new.Cells.Item(3, 6).Value = 'Alabama'; $new.Cells.Item(3, 7).Value = '6''9"'; $new.Cells.Item(3, 8).Value = '200 lbs'; $new.Cells.Item(3, 9).Value = 37582
    $new.Cells.Item(4, 1).Value = 3; $new.Cells.Item(4, 2).Value = 2; $new.Cells.Item(4, 3).Value = 'Scoot Henderson'; $new.Cells.Item(4, 4).Value = 'Guard'; $new.Cells.Item(4, 5).Value = 'Guard'; $new.Cells.Item(4, 6).Value = 'G League Ignite'; $new.Cells.Item(4, 7).Value = '6''2"'; $new.Cells.Item(4, 8).Value = '195 lbs'; $new.Cells.Item(4, 9).Value = 38020
    $new.Cells.Item(5, 1).Value = 4; $new.Cells.Item(5, 2).Value = 3; $new.Cells.Item(5, 3).Value = 'Jarace Walker'; $new.Cells.Item(5, 4).Value = 'Forward'; $new.Cells.Item(5, 5).Value = 'Forward'; $new.Cells.Item(5, 6).Value = 'Houston'; $new.Cells.Item(5, 7).Value = '6''7"'; $new.Cells.Item(5, 8).Value = '249 lbs'; $new.Cells.Item(5, 9).Value = 37868
    $new.Cells.Item(6, 1).Value = 5; $new.Cells.Item(6, 2).Value = 3; $new.Cells.Item(6, 3).Value = 'Cam Whitmore'; $new.Cells.Item(6, 4).Value = 'Wing'; $new.Cells.Item(6, 5).Value = 'Forward'; $new.Cells.Item(6, 6).Value = 'Villanova'; $new.Cells.Item(6, 7).Value = '6''6"'; $new.Cells.Item(6, 8).Value = '235 lbs'; $new.Cells.Item(6, 9).Value = 38176
    $new.Cells.Item(7, 1).Value = 6; $new.Cells.Item(7, 2).Value = 3; $new.Cells.Item(7, 3).Value = 'Ausar Thompson'; $new.Cells.Item(7, 4).Value = 'Wing'; $new.Cells.Item(7, 5).Value = 'Wing'; $new.Cells.Item(7, 6).Value = 'Overtime Elite'; $new.Cells.Item(7, 7).Value = '6''6"'; $new.Cells.Item(7, 8).Value = '218 lbs'; $new.Cells.Item(7, 9).Value = 37651
    $new.Cells.Item(8, 1).Value = 7; $new.Cells.Item(8, 2).Value = 3; $new.Cells.Item(8, 3).Value = 'Amen Thompson'; $new.Cells.Item(8, 4).Value = 'Guard'; $new.Cells.Item(8, 5).Value = 'Wing'; $new.Cells.Item(8, 6).Value = 'Overtime Elite'; $new.Cells.Item(8, 7).Value = '6''6"'; $new.Cells.Item(8, 8).Value = '214 lbs'; $new.Cells.Item(8, 9).Value = 37651
    $new.Cells.Item(9, 1).Value = 8; $new.Cells.Item(9, 2).Value = 3; $new.Cells.Item(9, 3).Value = 'Anthony Black'; $new.Cells.Item(9, 4).Value = 'Guard'; $new.Cells.Item(9, 5).Value = 'Guard'; $new.Cells.Item(9, 6).Value = 'Arkansas'; $new.Cells.Item(9, 7).Value = '6''6"'; $new.Cells.Item(9, 8).Value = '210 lbs'; $new.Cells.Item(9, 9).Value = 38006
    $new.Cells.Item(10, 1).Value = 9; $new.Cells.Item(10, 2).Value = 3; $new.Cells.Item(10, 3).Value = 'Cason Wallace'; $new.Cells.Item(10, 4).Value = 'Guard'; $new.Cells.Item(10, 5).Value = 'Guard'; $new.Cells.Item(10, 6).Value = 'Kentucky'; $new.Cells.Item(10, 7).Value = '6''3"'; $new.Cells.Item(10, 8).Value = '195 lbs'; $new.Cells.Item(10, 9).Value = 37932
    $new.Cells.Item(11, 1).Value = 10; $new.Cells.Item(11, 2).Value = 3; $new.Cells.Item(11, 3).Value = 'Taylor Hendricks'; $new.Cells.Item(11, 4).Value = 'Forward'; $new.Cells.Item(11, 5).Value = 'Forward'; $new.Cells.Item(11, 6).Value = 'UCF'; $new.Cells.Item(11, 7).Value = '6''8"'; $new.Cells.Item(11, 8).Value = '214 lbs'; $new.Cells.Item(11, 9).Value = 37947
    $new.Cells.Item(12, 1).Value = 11; $new.Cells.Item(12, 2).Value = 4; $new.Cells.Item(12, 3).Value = 'Jalen Hood-Schifino'; $new.Cells.Item(12, 4).Value = 'Guard'; $new.Cells.Item(12, 5).Value = 'Guard'; $new.Cells.Item(12, 6).Value = 'Indiana'; $new.Cells.Item(12, 7).Value = '6''4"'; $new.Cells.Item(12, 8).Value = '217 lbs'; $new.Cells.Item(12, 9).Value = 37791
    $new.Cells.Item(13, 1).Value = 12; $new.Cells.Item(13, 2).Value = 4; $new.Cells.Item(13, 3).Value = 'Kobe Bufkin'; $new.Cells.Item(13, 4).Value = 'Guard'; $new.Cells.Item(13, 5).Value = 'Guard'; $new.Cells.Item(13, 6).Value = 'Michigan'; $new.Cells.Item(13, 7).Value = '6''4"'; $new.Cells.Item(13, 8).Value = '187 lbs'; $new.Cells.Item(13, 9).Value = 37885
    $new.Cells.Item(14, 1).Value = 13; $new.Cells.Item(14, 2).Value = 4; $new.Cells.Item(14, 3).Value = 'Dariq Whitehead'; $new.Cells.Item(14, 4).Value = 'Wing'; $new.Cells.Item(14, 5).Value = 'Wing'; $new.Cells.Item(14, 6).Value = 'Duke'; $new.Cells.Item(14, 7).Value = '6''6"'; $new.Cells.Item(14, 8).Value = '217 lbs'; $new.Cells.Item(14, 9).Value = 38200
    $new.Cells.Item(15, 1).Value = 14; $new.Cells.Item(15, 2).Value = 4; $new.Cells.Item(15, 3).Value = 'Maxwell Lewis'; $new.Cells.Item(15, 4).Value = 'Wing'; $new.Cells.Item(15, 5).Value = 'Wing'; $new.Cells.Item(15, 6).Value = 'Pepperdine'; $new.Cells.Item(15, 7).Value = '6''6"'; $new.Cells.Item(15, 8).Value = '207 lbs'; $new.Cells.Item(15, 9).Value = 37464
    $new.Cells.Item(16, 1).Value = 15; $new.Cells.Item(16, 2).Value = 4; $new.Cells.Item(16, 3).Value = 'Keyonte George'; $new.Cells.Item(16, 4).Value = 'Guard'; $new.Cells.Item(16, 5).Value = 'Guard'; $new.Cells.Item(16, 6).Value = 'Baylor'; $new.Cells.Item(16, 7).Value = '6''4"'; $new.Cells.Item(16, 8).Value = '185 lbs'; $new.Cells.Item(16, 9).Value = 37933
    $new.Cells.Item(17, 1).Value = 16; $new.Cells.Item(17, 2).Value = 4; $new.Cells.Item(17, 3).Value = 'GG Jackson'; $new.Cells.Item(17, 4).Value = 'Wing'; $new.Cells.Item(17, 5).Value = 'Forward'; $new.Cells.Item(17, 6).Value = 'South Carolina'; $new.Cells.Item(17, 7).Value = '6''8"'; $new.Cells.Item(17, 8).Value = '214 lbs'; $new.Cells.Item(17, 9).Value = 38338
    $new.Cells.Item(18, 1).Value = 17; $new.Cells.Item(18, 2).Value = 4; $new.Cells.Item(18, 3).Value = 'Leonard Miller'; $new.Cells.Item(18, 4).Value = 'Forward'; $new.Cells.Item(18, 5).Value = 'Forward'; $new.Cells.Item(18, 6).Value = 'G League Ignite'; $new.Cells.Item(18, 7).Value = '6''9"'; $new.Cells.Item(18, 8).Value = '213 lbs'; $new.Cells.Item(18, 9).Value = 37951
    $new.Cells.Item(19, 1).Value = 18; $new.Cells.Item(19, 2).Value = 4; $new.Cells.Item(19, 3).Value = 'Jett Howard'; $new.Cells.Item(19, 4).Value = 'Guard'; $new.Cells.Item(19, 5).Value = 'Wing'; $new.Cells.Item(19, 6).Value = 'Michigan'; $new.Cells.Item(19, 7).Value = '6''8"'; $new.Cells.Item(19, 8).Value = '215 lbs'; $new.Cells.Item(19, 9).Value = 37878
    $new.Cells.Item(20, 1).Value = 19; $new.Cells.Item(20, 2).Value = 4; $new.Cells.Item(20, 3).Value = 'Jordan Hawkins'; $new.Cells.Item(20, 4).Value = 'Guard'; $new.Cells.Item(20, 5).Value = 'Guard'; $new.Cells.Item(20, 6).Value = 'UConn'; $new.Cells.Item(20, 7).Value = '6''4"'; $new.Cells.Item(20, 8).Value = '186 lbs'; $new.Cells.Item(20, 9).Value = 37375
    $new.Cells.Item(21, 1).Value = 20; $new.Cells.Item(21, 2).Value = 4; $new.Cells.Item(21, 3).Value = 'Bilal Coulibaly'; $new.Cells.Item(21, 4).Value = 'Wing'; $new.Cells.Item(21, 5).Value = 'Wing'; $new.Cells.Item(21, 6).Value = 'Metropolitans92'; $new.Cells.Item(21, 7).Value = '6''7"'; $new.Cells.Item(21, 8).Value = '194 lbs'; $new.Cells.Item(21, 9).Value = 38194
    $new.Cells.Item(22, 1).Value = 21; $new.Cells.Item(22, 2).Value = 4; $new.Cells.Item(22, 3).Value = 'Gradey Dick'; $new.Cells.Item(22, 4).Value = 'Guard'; $new.Cells.Item(22, 5).Value = 'Wing'; $new.Cells.Item(22, 6).Value = 'Kansas'; $new.Cells.Item(22, 7).Value = '6''6"'; $new.Cells.Item(22, 8).Value = '204 lbs'; $new.Cells.Item(22, 9).Value = 37945
    $new.Cells.Item(23, 1).Value = 22; $new.Cells.Item(23, 2).Value = 4; $new.Cells.Item(23, 3).Value = 'Dereck Lively II'; $new.Cells.Item(23, 4).Value = 'Big'; $new.Cells.Item(23, 5).Value = 'Big'; $new.Cells.Item(23, 6).Value = 'Duke'; $new.Cells.Item(23, 7).Value = '7''1"'; $new.Cells.Item(23, 8).Value = '230 lbs'; $new.Cells.Item(23, 9).Value = 38029
    $new.Cells.Item(24, 1).Value = 23; $new.Cells.Item(24, 2).Value = 4; $new.Cells.Item(24, 3).Value = 'Sidy Cissoko'; $new.Cells.Item(24, 4).Value = 'Forward'; $new.Cells.Item(24, 5).Value = 'Forward'; $new.Cells.Item(24, 6).Value = 'G League Ignite'; $new.Cells.Item(24, 7).Value = '6''6"'; $new.Cells.Item(24, 8).Value = '224 lbs'; $new.Cells.Item(24, 9).Value = 38079
    $new.Cells.Item(25, 1).Value = 24; $new.Cells.Item(25, 2).Value = 4; $new.Cells.Item(25, 3).Value = 'Brice Sensabaugh'; $new.Cells.Item(25, 4).Value = 'Wing'; $new.Cells.Item(25, 5).Value = 'Wing'; $new.Cells.Item(25, 6).Value = 'Ohio State'; $new.Cells.Item(25, 7).Value = '6''6"'; $new.Cells.Item(25, 8).Value = '235 lbs'; $new.Cells.Item(25, 9).Value = 37924
    $new.Cells.Item(26, 1).Value = 25; $new.Cells.Item(26, 2).Value = 4; $new.Cells.Item(26, 3).Value = 'Nick Smith Jr.'; $new.Cells.Item(26, 4).Value = 'Guard'; $new.Cells.Item(26, 5).Value = 'Guard'; $new.Cells.Item(26, 6).Value = 'Arkansas'; $new.Cells.Item(26, 7).Value = '6''5"'; $new.Cells.Item(26, 8).Value = '185 lbs'; $new.Cells.Item(26, 9).Value = 38095
    $new.Cells.Item(27, 1).Value = 26; $new.Cells.Item(27, 2).Value = 5; $new.Cells.Item(27, 3).Value = 'Colby Jones'; $new.Cells.Item(27, 4).Value = 'Guard'; $new.Cells.Item(27, 5).Value = 'Wing'; $new.Cells.Item(27, 6).Value = 'Xavier'; $new.Cells.Item(27, 7).Value = '6''5"'; $new.Cells.Item(27, 8).Value = '199 lbs'; $new.Cells.Item(27, 9).Value = 37404
    $new.Cells.Item(28, 1).Value = 27; $new.Cells.Item(28, 2).Value = 5; $new.Cells.Item(28, 3).Value = 'Kris Murray'; $new.Cells.Item(28, 4).Value = 'Forward'; $new.Cells.Item(28, 5).Value = 'Forward'; $new.Cells.Item(28, 6).Value = 'Iowa State'; $new.Cells.Item(28, 7).Value = '6''8"'; $new.Cells.Item(28, 8).Value = '213 lbs'; $new.Cells.Item(28, 9).Value = 36757
    $new.Cells.Item(29, 1).Value = 28; $new.Cells.Item(29, 2).Value = 5; $new.Cells.Item(29, 3).Value = 'James Nnaji'; $new.Cells.Item(29, 4).Value = 'Big'; $new.Cells.Item(29, 5).Value = 'Big'; $new.Cells.Item(29, 6).Value = 'Barcelona (ACB)'; $new.Cells.Item(29, 7).Value = '6''11"'; $new.Cells.Item(29, 8).Value = '251 lbs'; $new.Cells.Item(29, 9).Value = 38213
    $new.Cells.Item(30, 1).Value = 29; $new.Cells.Item(30, 2).Value = 5; $new.Cells.Item(30, 3).Value = 'Andre Jackson'; $new.Cells.Item(30, 4).Value = 'Forward'; $new.Cells.Item(30, 5).Value = 'Forward'; $new.Cells.Item(30, 6).Value = 'UConn'; $new.Cells.Item(30, 7).Value = '6''6"'; $new.Cells.Item(30, 8).Value = '198 lbs'; $new.Cells.Item(30, 9).Value = 37208
    $new.Cells.Item(31, 1).Value = 30; $new.Cells.Item(31, 2).Value = 5; $new.Cells.Item(31, 3).Value = 'Brandin Podziemski'; $new.Cells.Item(31, 4).Value = 'Wing'; $new.Cells.Item(31, 5).Value = 'Wing'; $new.Cells.Item(31, 6).Value = 'Santa Clara'; $new.Cells.Item(31, 7).Value = '6''4"'; $new.Cells.Item(31, 8).Value = '204 lbs'; $new.Cells.Item(31, 9).Value = 37677
    $new.Cells.Item(32, 1).Value = 31; $new.Cells.Item(32, 2).Value = 5; $new.Cells.Item(32, 3).Value = 'Jaylen Clark'; $new.Cells.Item(32, 4).Value = 'Wing'; $new.Cells.Item(32, 5).Value = 'Wing'; $new.Cells.Item(32, 6).Value = 'UCLA'; $new.Cells.Item(32, 7).Value = '6''4"'; $new.Cells.Item(32, 8).Value = '204 lbs'; $new.Cells.Item(32, 9).Value = 37177
    $new.Cells.Item(33, 1).Value = 32; $new.Cells.Item(33, 2).Value = 5; $new.Cells.Item(33, 3).Value = 'Marcus Sasser'; $new.Cells.Item(33, 4).Value = 'Guard'; $new.Cells.Item(33, 5).Value = 'Guard'; $new.Cells.Item(33, 6).Value = 'Houston'; $new.Cells.Item(33, 7).Value = '6''1"'; $new.Cells.Item(33, 8).Value = '196 lbs'; $new.Cells.Item(33, 9).Value = 36790
    $new.Cells.Item(34, 1).Value = 33; $new.Cells.Item(34, 2).Value = 5; $new.Cells.Item(34, 3).Value = 'Terquavion Smith'; $new.Cells.Item(34, 4).Value = 'Guard'; $new.Cells.Item(34, 5).Value = 'Guard'; $new.Cells.Item(34, 6).Value = 'NC State'; $new.Cells.Item(34, 7).Value = '6''3"'; $new.Cells.Item(34, 8).Value = '163 lbs'; $new.Cells.Item(34, 9).Value = 37621
    $new.Cells.Item(35, 1).Value = 34; $new.Cells.Item(35, 2).Value = 5; $new.Cells.Item(35, 3).Value = 'Noah Clowney'; $new.Cells.Item(35, 4).Value = 'Forward'; $new.Cells.Item(35, 5).Value = 'Forward'; $new.Cells.Item(35, 6).Value = 'Alabama'; $new.Cells.Item(35, 7).Value = '6''10"'; $new.Cells.Item(35, 8).Value = '210 lbs'; $new.Cells.Item(35, 9).Value = 38182
    $new.Cells.Item(36, 1).Value = 35; $new.Cells.Item(36, 2).Value = 5; $new.Cells.Item(36, 3).Value = 'Ben Sheppard'; $new.Cells.Item(36, 4).Value = 'Guard'; $new.Cells.Item(36, 5).Value = 'Wing'; $new.Cells.Item(36, 6).Value = 'Belmont'; $new.Cells.Item(36, 7).Value = '6''5"'; $new.Cells.Item(36, 8).Value = '195 lbs'; $new.Cells.Item(36, 9).Value = 36907
    $new.Cells.Item(37, 1).Value = 36; $new.Cells.Item(37, 2).Value = 6; $new.Cells.Item(37, 3).Value = 'Rayan Rupert'; $new.Cells.Item(37, 4).Value = 'Wing'; $new.Cells.Item(37, 5).Value = 'Wing'; $new.Cells.Item(37, 6).Value = 'NZ Breakers'; $new.Cells.Item(37, 7).Value = '6''6"'; $new.Cells.Item(37, 8).Value = '193 lbs'; $new.Cells.Item(37, 9).Value = 38138
    $new.Cells.Item(38, 1).Value = 37; $new.Cells.Item(38, 2).Value = 6; $new.Cells.Item(38, 3).Value = 'Amari Bailey'; $new.Cells.Item(38, 4).Value = 'Guard'; $new.Cells.Item(38, 5).Value = 'Guard'; $new.Cells.Item(38, 6).Value = 'UCLA'; $new.Cells.Item(38, 7).Value = '6''3"'; $new.Cells.Item(38, 8).Value = '191 lbs'; $new.Cells.Item(38, 9).Value = 38034
    $new.Cells.Item(39, 1).Value = 38; $new.Cells.Item(39, 2).Value = 6; $new.Cells.Item(39, 3).Value = 'Omari Moore'; $new.Cells.Item(39, 4).Value = 'Guard'; $new.Cells.Item(39, 5).Value = 'Guard'; $new.Cells.Item(39, 6).Value = 'San Jose State'; $new.Cells.Item(39, 7).Value = '6''5"'; $new.Cells.Item(39, 8).Value = '189 lbs'; $new.Cells.Item(39, 9).Value = 36787
    $new.Cells.Item(40, 1).Value = 39; $new.Cells.Item(40, 2).Value = 6; $new.Cells.Item(40, 3).Value = 'Ricky Council IV'; $new.Cells.Item(40, 4).Value = 'Wing'; $new.Cells.Item(40, 5).Value = 'Wing'; $new.Cells.Item(40, 6).Value = 'Arkansas'; $new.Cells.Item(40, 7).Value = '6''5"'; $new.Cells.Item(40, 8).Value = '208 lbs'; $new.Cells.Item(40, 9).Value = 37106
    $new.Cells.Item(41, 1).Value = 40; $new.Cells.Item(41, 2).Value = 6; $new.Cells.Item(41, 3).Value = 'Trayce Jackson-Davis'; $new.Cells.Item(41, 4).Value = 'Big'; $new.Cells.Item(41, 5).Value = 'Big'; $new.Cells.Item(41, 6).Value = 'Indiana'; $new.Cells.Item(41, 7).Value = '6''8"'; $new.Cells.Item(41, 8).Value = '240 lbs'; $new.Cells.Item(41, 9).Value = 36578
    $new.Cells.Item(42, 1).Value = 41; $new.Cells.Item(42, 2).Value = 6; $new.Cells.Item(42, 3).Value = 'Mike Miles Jr.'; $new.Cells.Item(42, 4).Value = 'Guard'; $new.Cells.Item(42, 5).Value = 'Guard'; $new.Cells.Item(42, 6).Value = 'TCU'; $new.Cells.Item(42, 7).Value = '6''1"'; $new.Cells.Item(42, 8).Value = '205 lbs'; $new.Cells.Item(42, 9).Value = 37492
    $new.Cells.Item(43, 1).Value = 42; $new.Cells.Item(43, 2).Value = 6; $new.Cells.Item(43, 3).Value = 'Jalen Wilson'; $new.Cells.Item(43, 4).Value = 'Wing'; $new.Cells.Item(43, 5).Value = 'Forward'; $new.Cells.Item(43, 6).Value = 'Kansas'; $new.Cells.Item(43, 7).Value = '6''6"'; $new.Cells.Item(43, 8).Value = '230 lbs'; $new.Cells.Item(43, 9).Value = 36834
    $new.Cells.Item(44, 1).Value = 43; $new.Cells.Item(44, 2).Value = 6; $new.Cells.Item(44, 3).Value = 'Jaime Jaquez Jr.'; $new.Cells.Item(44, 4).Value = 'Forward'; $new.Cells.Item(44, 5).Value = 'Forward'; $new.Cells.Item(44, 6).Value = 'UCLA'; $new.Cells.Item(44, 7).Value = '6''6"'; $new.Cells.Item(44, 8).Value = '226 lbs'; $new.Cells.Item(44, 9).Value = 36940
    $new.Cells.Item(45, 1).Value = 44; $new.Cells.Item(45, 2).Value = 6; $new.Cells.Item(45, 3).Value = 'Olivier-Maxence Prosper'; $new.Cells.Item(45, 4).Value = 'Forward'; $new.Cells.Item(45, 5).Value = 'Forward'; $new.Cells.Item(45, 6).Value = 'Marquette'; $new.Cells.Item(45, 7).Value = '6''7"'; $new.Cells.Item(45, 8).Value = '212 lbs'; $new.Cells.Item(45, 9).Value = 37440

# The new sheet becomes the active/selected tab, with the same cell
# selection the author left it on.
$new.Activate()
$new.Range("E8").Select()

# The original board is no longer the active tab; scroll it over so column D
# is the leftmost visible column.
$orig = $wb.Worksheets.Item("Adam-Bushman")
$orig.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$new.Activate()
